# Highlight (yellow) the six pieces of requirement text that were marked
# up in the commit, leaving the trailing double-space at the end of each
# paragraph un-highlighted (matches the XML diff: Word splits the run at
# the highlight boundary).
#
# wdYellow = 7 (HighlightColorIndex enum)
$wdYellow = 7
$wdFindStop = 0

$d = $word.ActiveDocument

# Use a single range that we keep collapsing to the end of the previous
# match and re-searching forward from there, so the two occurrences of
# the word "troquer" are each matched in document order (first one is
# left untouched, second one - the underlined hyperlink placeholder -
# gets the highlight).
$rng = $d.Content

function Find-And-Highlight($range, $text) {
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($text, $true, $false, $false, $false, $false, $true, $wdFindStop, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $text"
    }
    $range.Font.HighlightColorIndex = $wdYellow
    $range.Collapse(0)
}

Find-And-Highlight $rng "Lorsqu’un utilisateur est authentifié, il est possible d’ajouter, modifier et supprimer un objet. Il ne peut modifier ou supprimer que ceux qui lui appartiennent."

Find-And-Highlight $rng "Un utilisateur authentifié peut troquer un objet avec un objet de son choix. Ajoutez un lien dans la liste d’objets (visible seulement aux personnes identifiées) à côté de chaque objet dont le texte est "

Find-And-Highlight $rng "troquer"

Find-And-Highlight $rng ". Ce lien ouvre une page qui contient une liste déroulante des objets de l’utilisateur authentifié. Il en choisit un et le propriétaire de chacun des objets est interverti."

Find-And-Highlight $rng "Seul un compte administrateur ou l’utilisateur qui a créé un objet peuvent l’effacer ou le modifier."

Find-And-Highlight $rng "Un compte non authentifié qui essaie d’aller ailleurs que sur la page d’accueil, la liste d’objet, la page d’authentification ou la page de détail d’un objet doit être bloqué."

Write-Output "done"
